# Added confidence level after each extraction
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D8 was empty -> set to "Expansion Tank"
$ws.Range("D8").Value = "Expansion Tank"

# Normalize FLUID column text from all-caps to capitalized, rows 8-15
$ws.Range("G8").Value = "Condensate"
$ws.Range("G9").Value = "Condensate"
$ws.Range("G10").Value = "Condensate"
$ws.Range("G11").Value = "Condensate"
$ws.Range("G12").Value = "Condensate"
$ws.Range("G13").Value = "Condensate"
$ws.Range("G14").Value = "Condensate"
$ws.Range("G15").Value = "Condensate"

# GRADE column: strip "Gr." prefix (force text so "70" stays a string, not a number)
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "70"
$ws.Range("J10").Value = "B"
$ws.Range("J12").NumberFormat = "@"
$ws.Range("J12").Value = "70"

# I13/J13: split "S275JR" into SPEC "S275" and GRADE "JR"
$ws.Range("I13").Value = "S275"
$ws.Range("J13").Value = "JR"

# DESIGN / OPERATING TEMP columns: add space before "C", rows 8-15
$ws.Range("L8").Value = "200° C"
$ws.Range("N8").Value = "185° C"
$ws.Range("L9").Value = "200° C"
$ws.Range("N9").Value = "185° C"
$ws.Range("L10").Value = "200° C"
$ws.Range("N10").Value = "185° C"
$ws.Range("L11").Value = "200° C"
$ws.Range("N11").Value = "185° C"
$ws.Range("L12").Value = "200° C"
$ws.Range("N12").Value = "185° C"
$ws.Range("L13").Value = "200° C"
$ws.Range("N13").Value = "185° C"
$ws.Range("L14").Value = "200° C"
$ws.Range("N14").Value = "185° C"
$ws.Range("L15").Value = "200° C"
$ws.Range("N15").Value = "185° C"

# PARTS column renames
$ws.Range("E14").Value = "Bolt & Nut (Pressure Retaining)"
$ws.Range("E15").Value = "Bolt & Nut (External Fittings)"
